# edit.ps1 -- applies the "updating plan for results" edits:
#   1. Insert a new outline bullet "Add other citations for using ratios."
#      (one indent level deeper than its neighbour) right before the
#      "Generally, increased residence time..." bullet.
#   2. Insert the phrase "relative to trophic state " into the
#      "How does nutrient limitation ... vary spatially..." question.
#   3. Add a reviewer comment ("Move to Q1. and Q1. to Q3.", by Linnea Ann
#      Rock) anchored on the "What are the trends ... scale?" question.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. New bullet: "Add other citations for using ratios."
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Generally, increased residence time correlates", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertParagraphBefore()

# Re-acquire the (now empty) paragraph that was just created and fill it in.
$newPara = $d.Range($rng.Start, $rng.Start).Paragraphs(1)
$newParaRange = $newPara.Range
$fillRange = $d.Range($newParaRange.Start, $newParaRange.End - 1)
$fillRange.InsertAfter("Add other citations for using ratios.  ")
# Bump it one level deeper than the surrounding ilvl=2 bullets (-> ilvl=3).
$fillRange.ListFormat.ListLevelNumber = 4

# ---------------------------------------------------------------------------
# 2. "relative to trophic state " inserted into the Q1 question.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute(" vary spatially and temporally across the US? ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertion = $d.Range($rng2.Start + 1, $rng2.Start + 1)
$insertion.InsertBefore("relative to trophic state ")

# ---------------------------------------------------------------------------
# 3. Reviewer comment on the Q3 question.
# ---------------------------------------------------------------------------
$word.UserName = "Linnea Ann Rock"
$word.UserInitials = "LAR"

$qStart = $d.Content
$null = $qStart.Find.Execute("What are the trends", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $qStart.Start

$qEnd = $d.Content
$null = $qEnd.Find.Execute("ross ecoregional and the national scale", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $qEnd.End + 1

$commentRange = $d.Range($startPos, $endPos)
$comment = $d.Comments.Add($commentRange, "Move to Q1. and Q1. to Q3.")

Write-Host "Edits applied."
